$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")

# Insert a new row above the existing row 22 ("Bekannte Sensorausfälle:") to host
# the new "Störungen:" remark about the West platform counters.
$ws.Rows.Item(22).Insert()

$ws.Range("A22").Value = $ws.Range("A23").Value
$ws.Range("B22").Value = "Störungen:"
$ws.Range("C22").Value = "'- Zählwerte ab 1.1.2023 auf dem **Perron ""West""** sind **nicht korrekt**. Bitte verwenden Sie die absoluten Werte derzeit nicht, bis das Problem behoben ist."

# Match styling used by the other "bemerkung" rows in this block.
$ws.Range("A22").Style = $ws.Range("A23").Style
$ws.Range("B22").Style = $ws.Range("B23").Style
$ws.Range("D22").Style = $ws.Range("D23").Style
$ws.Range("E22").Style = $ws.Range("E23").Style
$ws.Range("F22").Style = $ws.Range("F23").Style

$ws.Rows.Item(22).RowHeight = 82.8
